$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet: handback status text updated for both zh-cn / de-de columns ---
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# Widen the zh-cn / de-de status columns on the Overview sheet
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# --- zh-cn sheet updates ---
# Status column (C) gets wider
$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
# Error Detail column (P) gets narrower
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333332
# Latest Handback DateTime (K2) refreshed
$wsZhCn.Range("K2").Value = "2016-09-01 02:58:54"
# Error Detail (P2) cleared - handback is now in sync, no error
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet updates ---
# Status column (C) gets wider
$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
# Error Detail column (P) gets narrower
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333332
# Latest Handback DateTime (K2) refreshed
$wsDeDe.Range("K2").Value = "2016-09-01 02:59:06"
# Error Detail (P2) cleared - handback is now in sync, no error
$wsDeDe.Range("P2").Value = ""
